$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.206.86'
$ws.Range("E2").Value = '  -3.33%  '

$ws.Range("D3").Value = '1.924.47'
$ws.Range("E3").Value = '  -2.89%  '

$ws.Range("D4").Value = '''0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '''245.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.77%  '

$ws.Range("D6").Value = '''0.7192'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -10.55%  '

$ws.Range("D7").Value = '''0.9993'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = '''0.3241'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.36%  '

$ws.Range("D9").Value = '''26.43'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.10%  '

$ws.Range("D10").Value = '''0.06836'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").Value = '''0.7973'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.46%  '

$ws.Range("D12").Value = '''0.07924'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.21%  '

$ws.Range("D13").Value = '1.921.62'
$ws.Range("E13").Value = '  -3.00%  '

$ws.Range("D14").Value = '''5.386'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").Value = '''94.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.61%  '

$ws.Range("D16").Value = '''14.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.36%  '

$ws.Range("D17").Value = '''260.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.58%  '

$ws.Range("D18").Value = '30.212.99'
$ws.Range("E18").Value = '  -3.31%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '''0.000007924'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''5.806'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").Value = '2.174.25'
$ws.Range("E21").Value = '  -3.14%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = '''0.9992'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").Value = '''6.850'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("D25").Value = '''9.640'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("D26").Value = '''160.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.30%  '

$ws.Range("D27").Value = '''0.1333'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.76%  '

$ws.Range("D28").Value = '''18.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.57%  '

$ws.Range("D29").Value = '''2.259'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.55%  '

$ws.Range("D30").Value = '''1.354'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("D31").Value = '''1.542'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.36%  '

$ws.Range("D32").Value = '''4.411'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.66%  '

$ws.Range("D33").Value = '''4.182'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.38%  '

$ws.Range("D34").Value = '''0.05047'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.17%  '

$ws.Range("E35").Value = '  -1.80%  '

$ws.Range("D36").Value = '''0.7370'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("D37").Value = '''2.737'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.89%  '

$ws.Range("D38").Value = '''0.01933'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.77%  '

$ws.Range("E39").Value = '  -3.57%  '

$ws.Range("D40").Value = '''80.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.53%  '

$ws.Range("D41").Value = '''6.519'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.23%  '

$ws.Range("D42").Value = '''0.4442'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.85%  '

$ws.Range("D43").Value = '''2.004'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.66%  '

$ws.Range("D44").Value = '''0.9994'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("D45").Value = '''0.8307'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.70%  '

$ws.Range("D46").Value = '''102.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.04%  '

$ws.Range("D47").Value = '''9.674'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.96%  '

$ws.Range("E48").Value = '  -2.94%  '

$ws.Range("D49").Value = '''36.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.95%  '

$ws.Range("E50").Value = '  -4.37%  '

$ws.Range("E51").Value = '  +2.43%  '
